$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 287414
$ws.Range("D2").Value = 366965635
$ws.Range("C3").Value = 235
$ws.Range("D3").Value = 280129
$ws.Range("C10").Value = 108958
$ws.Range("D10").Value = 159788908
$ws.Range("C12").Value = 53806
$ws.Range("D12").Value = 77718058
$ws.Range("C16").Value = 3600
$ws.Range("D16").Value = 5115962
$ws.Range("C20").Value = 5375
$ws.Range("D20").Value = 7515177
$ws.Range("C22").Value = 70672
$ws.Range("D22").Value = 88554377
$ws.Range("C28").Value = 30434
$ws.Range("D28").Value = 44576992
$ws.Range("C30").Value = 10513
$ws.Range("D30").Value = 15146121
$ws.Range("C33").Value = 1411
$ws.Range("D33").Value = 1983144
$ws.Range("C35").Value = 1508
$ws.Range("D35").Value = 2126579
$ws.Range("C36").Value = 89323
$ws.Range("D36").Value = 113005450
$ws.Range("C42").Value = 861
$ws.Range("D42").Value = 1269461
$ws.Range("C44").Value = 41730
$ws.Range("D44").Value = 61224582
$ws.Range("C46").Value = 8360
$ws.Range("D46").Value = 12007299
$ws.Range("C48").Value = 1274
$ws.Range("D48").Value = 1767466
$ws.Range("C51").Value = 1894
$ws.Range("D51").Value = 2632984
$ws.Range("C52").Value = 62614
$ws.Range("D52").Value = 78753803
$ws.Range("C56").Value = 362
$ws.Range("D56").Value = 532539
$ws.Range("C58").Value = 26270
$ws.Range("D58").Value = 38543791
$ws.Range("C61").Value = 10122
$ws.Range("D61").Value = 14642985
$ws.Range("C63").Value = 1249
$ws.Range("D63").Value = 1743474
$ws.Range("C66").Value = 1218
$ws.Range("D66").Value = 1704455
$ws.Range("C68").Value = 18198
$ws.Range("D68").Value = 23806317
$ws.Range("C72").Value = 6647
$ws.Range("D72").Value = 9724521
$ws.Range("C74").Value = 4522
$ws.Range("D74").Value = 6569350
$ws.Range("C77").Value = 127941
$ws.Range("D77").Value = 159964918
$ws.Range("C83").Value = 59405
$ws.Range("D83").Value = 87156247
$ws.Range("C86").Value = 27341
$ws.Range("D86").Value = 39566345
$ws.Range("C88").Value = 2468
$ws.Range("D88").Value = 3557525
$ws.Range("C89").Value = 2308
$ws.Range("D89").Value = 3254673
$ws.Range("C90").Value = 26229
$ws.Range("D90").Value = 35566253
$ws.Range("C94").Value = 6698
$ws.Range("D94").Value = 9870874
$ws.Range("C96").Value = 5897
$ws.Range("D96").Value = 8541412
$ws.Range("C98").Value = 418
$ws.Range("D98").Value = 592049
$ws.Range("C100").Value = 6289
$ws.Range("D100").Value = 8706226
$ws.Range("C102").Value = 1574
$ws.Range("D102").Value = 2316907
$ws.Range("C104").Value = 2112
$ws.Range("D104").Value = 3075017
$ws.Range("C106").Value = 79
$ws.Range("D106").Value = 113870
$ws.Range("C108").Value = 128583
$ws.Range("D108").Value = 159185832
$ws.Range("C114").Value = 49458
$ws.Range("D114").Value = 72563670
$ws.Range("C116").Value = 24579
$ws.Range("D116").Value = 35613786
$ws.Range("C120").Value = 1875
$ws.Range("D120").Value = 2631977
$ws.Range("C122").Value = 414664
$ws.Range("D122").Value = 545008503
$ws.Range("C128").Value = 23
$ws.Range("D128").Value = 29273
$ws.Range("C129").Value = 186079
$ws.Range("D129").Value = 273762262
$ws.Range("C132").Value = 155942
$ws.Range("D132").Value = 226675823
$ws.Range("C135").Value = 2390
$ws.Range("D135").Value = 3349820
$ws.Range("C137").Value = 5002
$ws.Range("D137").Value = 7050337
$ws.Range("C140").Value = 38734
$ws.Range("D140").Value = 51871227
$ws.Range("C146").Value = 12817
$ws.Range("D146").Value = 18829297
$ws.Range("C147").Value = 3339
$ws.Range("D147").Value = 4821446
$ws.Range("C152").Value = 326
$ws.Range("D152").Value = 457159
$ws.Range("C153").Value = 14912
$ws.Range("D153").Value = 19779305
$ws.Range("C157").Value = 6304
$ws.Range("D157").Value = 9185514
$ws.Range("C159").Value = 4208
$ws.Range("D159").Value = 6070161
$ws.Range("C162").Value = 205
$ws.Range("D162").Value = 292501
$ws.Range("C164").Value = 10918
$ws.Range("D164").Value = 15780204
$ws.Range("C165").Value = 1392
$ws.Range("D165").Value = 2069578
$ws.Range("C169").Value = 79971
$ws.Range("D169").Value = 100351402
$ws.Range("C176").Value = 31770
$ws.Range("D176").Value = 46615495
$ws.Range("C178").Value = 11887
$ws.Range("D178").Value = 17185234
$ws.Range("C180").Value = 1126
$ws.Range("D180").Value = 1574620
$ws.Range("C182").Value = 1362
$ws.Range("D182").Value = 1912960
$ws.Range("C184").Value = 215771
$ws.Range("D184").Value = 268784470
$ws.Range("C192").Value = 81103
$ws.Range("D192").Value = 118937076
$ws.Range("C193").Value = 87
$ws.Range("D193").Value = 126127
$ws.Range("C195").Value = 30224
$ws.Range("D195").Value = 43511345
$ws.Range("C198").Value = 4601
$ws.Range("D198").Value = 6558146
$ws.Range("C201").Value = 4007
$ws.Range("D201").Value = 5555592
$ws.Range("C204").Value = 236145
$ws.Range("D204").Value = 292701739
$ws.Range("C206").Value = 230
$ws.Range("D206").Value = 328411
$ws.Range("C211").Value = 573
$ws.Range("D211").Value = 836127
$ws.Range("C213").Value = 88201
$ws.Range("D213").Value = 129106356
$ws.Range("C216").Value = 46803
$ws.Range("D216").Value = 67693849
$ws.Range("C219").Value = 4226
$ws.Range("D219").Value = 5932545
$ws.Range("C222").Value = 4629
$ws.Range("D222").Value = 6391417
$ws.Range("C225").Value = 96478
$ws.Range("D225").Value = 121180419
$ws.Range("C232").Value = 46452
$ws.Range("D232").Value = 68087650
$ws.Range("C234").Value = 11242
$ws.Range("D234").Value = 16170143
$ws.Range("C236").Value = 1768
$ws.Range("D236").Value = 2535733
$ws.Range("C238").Value = 2130
$ws.Range("D238").Value = 2965963
$ws.Range("C239").Value = 231227
$ws.Range("D239").Value = 292265639
$ws.Range("C240").Value = 153
$ws.Range("D240").Value = 189301
$ws.Range("C241").Value = 229
$ws.Range("D241").Value = 329842
$ws.Range("C247").Value = 88715
$ws.Range("D247").Value = 130082360
$ws.Range("C250").Value = 58840
$ws.Range("D250").Value = 85323944
$ws.Range("C252").Value = 2183
$ws.Range("D252").Value = 3080077
$ws.Range("C255").Value = 3804
$ws.Range("D255").Value = 5336163
